$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "passenger"
$ws.Range("E1").Value = "class"
$ws.Range("D2").Value = "'2"
$ws.Range("E2").Value = "Economy"
$ws.Range("D3").Value = "'2"
$ws.Range("E3").Value = "Economy"

$ws.Range("A2").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("E2:E3").PasteSpecial(-4122)

$ws.Range("M15").Select()
